$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 100000456
$ws.Range("I28").Value = 142857520
$ws.Range("J28").Value = 626.3333
$ws.Range("K28").Value = 142857520
$ws.Range("L28").Value = 626.3333
$ws.Range("M28").Value = -142857035
$ws.Range("N28").Value = -1596.3333

$ws.Range("H51").Value = 3999.8333
$ws.Range("J51").Value = 3999.7778
$ws.Range("L51").Value = 3999.7778
$ws.Range("N51").Value = -4967.7778

$ws.Range("H57").Value = 134270.67
$ws.Range("J57").Value = 134270.67
$ws.Range("L57").Value = 402812.01
$ws.Range("N57").Value = -403810.01

$ws.Range("H95").Value = 27215
$ws.Range("J95").Value = 27215
$ws.Range("L95").Value = 27215
$ws.Range("N95").Value = -32707

$ws.Range("H127").Value = 24973.75
$ws.Range("I127").Value = 2995.3333
$ws.Range("K127").Value = 8985.999899999999
$ws.Range("M127").Value = -4025.999899999999

$ws.Range("H129").Value = 2435.25
$ws.Range("I129").Value = 2435.25
$ws.Range("K129").Value = 7305.75
$ws.Range("M129").Value = -2305.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 3829.75
$ws.Range("J8").Value = 3439.6667
$ws.Range("L8").Value = 3439.6667
$ws.Range("N8").Value = -3727.6667

$ws.Range("H47").Value = 39333
$ws.Range("J47").Value = 39333
$ws.Range("L47").Value = 39333
$ws.Range("N47").Value = -40783

$ws.Range("H96").Value = 27360.54
$ws.Range("J96").Value = 27360.54
$ws.Range("L96").Value = 27360.54
$ws.Range("N96").Value = -32852.54

$ws.Range("H102").Value = 6365.5713
$ws.Range("I102").Value = 5355.846
$ws.Range("K102").Value = 5355.846
$ws.Range("M102").Value = -3733.846

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H49").Value = 16599.666
$ws.Range("J49").Value = 16599.666
$ws.Range("L49").Value = 16599.666
$ws.Range("N49").Value = -17077.666

$ws.Range("H80").Value = 897.5294
$ws.Range("J80").Value = 884.3333
$ws.Range("L80").Value = 884.3333
$ws.Range("N80").Value = -2880.3333

$ws.Range("H83").Value = 897.5294
$ws.Range("J83").Value = 884.3333
$ws.Range("L83").Value = 4421.6665
$ws.Range("N83").Value = -14405.6665

$ws.Range("H86").Value = 2178.1667
$ws.Range("I86").Value = 2170.5
$ws.Range("J86").Value = 2187.75
$ws.Range("K86").Value = 2170.5
$ws.Range("L86").Value = 2187.75
$ws.Range("M86").Value = -1047.5
$ws.Range("N86").Value = -4433.75

$ws.Range("H89").Value = 2178.1667
$ws.Range("I89").Value = 2170.5
$ws.Range("J89").Value = 2187.75
$ws.Range("K89").Value = 10852.5
$ws.Range("L89").Value = 10938.75
$ws.Range("M89").Value = -5236.5
$ws.Range("N89").Value = -22170.75

$ws.Range("H94").Value = 1454.5834
$ws.Range("I94").Value = 664.5
$ws.Range("K94").Value = 664.5
$ws.Range("M94").Value = -213.5

$ws.Range("H122").Value = 190000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 190000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 190000
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -199800

$ws.Range("H134").Value = 3078.6829
$ws.Range("I134").Value = 2283.9722
$ws.Range("J134").Value = 8800.6
$ws.Range("K134").Value = 6851.9166
$ws.Range("L134").Value = 26401.8
$ws.Range("M134").Value = -4316.9166
$ws.Range("N134").Value = -31471.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2619.6365
$ws.Range("I16").Value = 2053.7778
$ws.Range("K16").Value = 2053.7778
$ws.Range("M16").Value = -1766.7778

$ws.Range("H42").Value = 3287.4443
$ws.Range("I42").Value = 3287.4443
$ws.Range("K42").Value = 3287.4443
$ws.Range("M42").Value = -2694.4443

$ws.Range("H48").Value = 65518
$ws.Range("I48").Value = 51046
$ws.Range("K48").Value = 51046
$ws.Range("M48").Value = -50570

$ws.Range("H64").Value = 79999
$ws.Range("J64").Value = 79999
$ws.Range("L64").Value = 79999
$ws.Range("N64").Value = -80495

$ws.Range("H67").Value = 79999
$ws.Range("J67").Value = 79999
$ws.Range("L67").Value = 79999
$ws.Range("N67").Value = -81715

$ws.Range("H105").Value = 2087.6
$ws.Range("I105").Value = 2131.5454
$ws.Range("K105").Value = 2131.5454
$ws.Range("M105").Value = -384.5454

$ws.Range("H113").Value = 2619.6365
$ws.Range("I113").Value = 2053.7778
$ws.Range("K113").Value = 2053.7778
$ws.Range("M113").Value = 116.2222000000002

$ws.Range("H132").Value = 9934.641
$ws.Range("I132").Value = 11540.75
$ws.Range("J132").Value = 2592.4285
$ws.Range("K132").Value = 34622.25
$ws.Range("L132").Value = 7777.2855
$ws.Range("M132").Value = -32092.25
$ws.Range("N132").Value = -12837.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 12233
$ws.Range("J57").Value = 10224.5
$ws.Range("L57").Value = 30673.5
$ws.Range("N57").Value = -31791.5

$ws.Range("H117").Value = 3596.1428
$ws.Range("I117").Value = 4861
$ws.Range("K117").Value = 14583
$ws.Range("M117").Value = -11141

$ws.Range("H122").Value = 1647.8334
$ws.Range("I122").Value = 1353
$ws.Range("K122").Value = 12177
$ws.Range("M122").Value = -9727

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2220
$ws.Range("I70").Value = 2364.4
$ws.Range("J70").Value = 1601.1428
$ws.Range("K70").Value = 2364.4
$ws.Range("L70").Value = 1601.1428
$ws.Range("M70").Value = -2094.4
$ws.Range("N70").Value = -2141.1428

$ws.Range("H73").Value = 2220
$ws.Range("I73").Value = 2364.4
$ws.Range("J73").Value = 1601.1428
$ws.Range("K73").Value = 2364.4
$ws.Range("L73").Value = 1601.1428
$ws.Range("M73").Value = -1428.4
$ws.Range("N73").Value = -3473.1428

$ws.Range("H99").Value = 25614
$ws.Range("I99").Value = 28385.2
$ws.Range("J99").Value = 20995.334
$ws.Range("K99").Value = 28385.2
$ws.Range("L99").Value = 20995.334
$ws.Range("M99").Value = -26139.2
$ws.Range("N99").Value = -25487.334

$ws.Range("H113").Value = 4205.636
$ws.Range("I113").Value = 3677.25
$ws.Range("J113").Value = 4507.5713
$ws.Range("K113").Value = 3677.25
$ws.Range("L113").Value = 4507.5713
$ws.Range("M113").Value = -1507.25
$ws.Range("N113").Value = -8847.5713

$ws.Range("H122").Value = 4286.4287
$ws.Range("I122").Value = 6335.3335
$ws.Range("K122").Value = 19006.0005
$ws.Range("M122").Value = -16556.0005

$ws.Range("H126").Value = 5644.033
$ws.Range("I126").Value = 9974.25
$ws.Range("J126").Value = 4977.846
$ws.Range("K126").Value = 29922.75
$ws.Range("L126").Value = 14933.538
$ws.Range("M126").Value = -27452.75
$ws.Range("N126").Value = -19873.538

$ws.Range("H132").Value = 3716.111
$ws.Range("J132").Value = 4882.636
$ws.Range("L132").Value = 14647.908
$ws.Range("N132").Value = -19707.908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 17222
$ws.Range("J20").Value = 14444
$ws.Range("L20").Value = 14444
$ws.Range("N20").Value = -14896

$ws.Range("H22").Value = 940.5454999999999
$ws.Range("I22").Value = 953.6667
$ws.Range("K22").Value = 953.6667
$ws.Range("M22").Value = -658.6667

$ws.Range("H27").Value = 940.5454999999999
$ws.Range("I27").Value = 953.6667
$ws.Range("K27").Value = 953.6667
$ws.Range("M27").Value = -846.6667

$ws.Range("H46").Value = 2466.55
$ws.Range("J46").Value = 2363.6428
$ws.Range("L46").Value = 2363.6428
$ws.Range("N46").Value = -2739.6428

$ws.Range("H47").Value = 49000
$ws.Range("J47").Value = 49000
$ws.Range("L47").Value = 49000
$ws.Range("N47").Value = -49980

$ws.Range("H48").Value = 29994.5
$ws.Range("I48").Value = 29994.5
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 29994.5
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -29333.5
$ws.Range("N48").ClearContents()

$ws.Range("H52").Value = 49000
$ws.Range("J52").Value = 49000
$ws.Range("L52").Value = 49000
$ws.Range("N52").Value = -49466

$ws.Range("H55").Value = 800.4231
$ws.Range("I55").Value = 472.84616
$ws.Range("J55").Value = 1128
$ws.Range("K55").Value = 472.84616
$ws.Range("L55").Value = 1128
$ws.Range("M55").Value = -299.84616
$ws.Range("N55").Value = -1474

$ws.Range("H61").Value = 1279.5
$ws.Range("I61").Value = 1367.2
$ws.Range("K61").Value = 1367.2
$ws.Range("M61").Value = -1165.2

$ws.Range("H113").Value = 1279.5
$ws.Range("I113").Value = 1367.2
$ws.Range("K113").Value = 1367.2
$ws.Range("M113").Value = 802.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H107").Value = 43478476
$ws.Range("I107").Value = 202.44444
$ws.Range("K107").Value = 607.33332
$ws.Range("M107").Value = 1312.66668

$ws.Range("H122").Value = 3535.8
$ws.Range("I122").Value = 3633.2778
$ws.Range("K122").Value = 10899.8334
$ws.Range("M122").Value = -8449.8334
